# Updated cryptos list on Tue Jan 30 13:45:19 UTC 2024 with GitHub Actions
# Refreshes Price (col D) and Volume(1h) (col E) for each coin row on Sheet1,
# and (rows 39-40) swaps RenderToken/Kaspa back to their correct row order.
#
# Price values that look like a plain decimal number (e.g. "308.69") are written
# with a leading apostrophe so Excel stores them as literal text instead of
# silently parsing them into a float -- matching how the sheet already stores
# every Price/Volume cell as text (the column also holds values such as
# "43.237.94" that are not valid numbers at all).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '43.237.94'
$ws.Range('E2').Value = '  +2.81%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '2.302.72'
$ws.Range('E3').Value = '  +2.43%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  -0.01%  '

# Row 5: BNB
$ws.Range('D5').Value = '''308.69'
$ws.Range('E5').Value = '  +0.70%  '

# Row 6: Solana
$ws.Range('D6').Value = '''104.23'
$ws.Range('E6').Value = '  +8.12%  '

# Row 7: XRP
$ws.Range('E7').Value = '  +0.47%  '

# Row 9: Cardano
$ws.Range('D9').Value = '''0.520'
$ws.Range('E9').Value = '  +6.64%  '

# Row 10: Avalanche
$ws.Range('E10').Value = '  +3.15%  '

# Row 11: OKB
$ws.Range('D11').Value = '''52.10'
$ws.Range('E11').Value = '  +0.33%  '

# Row 12: Dogecoin
$ws.Range('D12').Value = '''0.0809'
$ws.Range('E12').Value = '  +0.06%  '

# Row 13: TRON
$ws.Range('E13').Value = '  -0.63%  '

# Row 14: Polkadot
$ws.Range('D14').Value = '''6.94'
$ws.Range('E14').Value = '  +2.51%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range('D15').Value = '2.660.63'
$ws.Range('E15').Value = '  +2.43%  '

# Row 16: Chainlink
$ws.Range('D16').Value = '''15.11'
$ws.Range('E16').Value = '  +4.59%  '

# Row 17: WrappedEther
$ws.Range('D17').Value = '2.302.83'
$ws.Range('E17').Value = '  +2.32%  '

# Row 18: Polygon
$ws.Range('D18').Value = '''0.803'
$ws.Range('E18').Value = '  +3.21%  '

# Row 19: WrappedBTC
$ws.Range('D19').Value = '43.191.96'
$ws.Range('E19').Value = '  +3.04%  '

# Row 20: InternetComputer(DFINITY)
$ws.Range('D20').Value = '''11.96'
$ws.Range('E20').Value = '  -1.50%  '

# Row 21: ShibaInu
$ws.Range('E21').Value = '  +2.22%  '

# Row 22: Uniswap
$ws.Range('E22').Value = '  +4.15%  '

# Row 23: Litecoin
$ws.Range('D23').Value = '''67.77'
$ws.Range('E23').Value = '  +0.96%  '

# Row 24: BitcoinCash
$ws.Range('D24').Value = '''240.07'
$ws.Range('E24').Value = '  +2.08%  '

# Row 25: ImmutableX
$ws.Range('E25').Value = '  +3.22%  '

# Row 26: PancakeSwap
$ws.Range('D26').Value = '''2.60'
$ws.Range('E26').Value = '  +1.38%  '

# Row 27: Dai
$ws.Range('E27').Value = '  +0.21%  '

# Row 28: EthereumClassic
$ws.Range('E28').Value = '  +5.73%  '

# Row 29: InjectiveProtocol
$ws.Range('E29').Value = '  -1.61%  '

# Row 30: Cosmos
$ws.Range('D30').Value = '''9.56'
$ws.Range('E30').Value = '  +0.89%  '

# Row 31: Toncoin
$ws.Range('D31').Value = '''2.11'
$ws.Range('E31').Value = '  -0.27%  '

# Row 32: Monero
$ws.Range('D32').Value = '''161.18'
$ws.Range('E32').Value = '  -2.36%  '

# Row 33: Filecoin
$ws.Range('D33').Value = '''5.24'
$ws.Range('E33').Value = '  +1.40%  '

# Row 34: FirstDigitalUSD
$ws.Range('D34').Value = '''0.999'
$ws.Range('E34').Value = '  -0.06%  '

# Row 35: Celestia
$ws.Range('D35').Value = '''18.21'
$ws.Range('E35').Value = '  +3.61%  '

# Row 36: WEMIXToken
$ws.Range('E36').Value = '  +6.79%  '

# Row 37: Hedera
$ws.Range('D37').Value = '''0.0735'
$ws.Range('E37').Value = '  +2.21%  '

# Row 38: LidoDAOToken
$ws.Range('D38').Value = '''3.00'
$ws.Range('E38').Value = '  -2.65%  '

# Row 39: RenderToken->Kaspa
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').Value = '''0.106'
$ws.Range('E39').Value = '  +3.74%  '

# Row 40: Kaspa->RenderToken
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').Value = '''4.46'
$ws.Range('E40').Value = '  +9.79%  '

# Row 41: ARBITRUM
$ws.Range('E41').Value = '  +3.53%  '

# Row 42: Stellar
$ws.Range('E42').Value = '  +0.50%  '

# Row 43: ApeXProtocol
$ws.Range('D43').Value = '''2.56'
$ws.Range('E43').Value = '  +18.98%  '

# Row 44: VeChain
$ws.Range('E44').Value = '  +2.89%  '

# Row 45: Maker
$ws.Range('D45').Value = '1.961.45'
$ws.Range('E45').Value = '  +0.90%  '

# Row 46: EnergySwap
$ws.Range('D46').Value = '''18.93'
$ws.Range('E46').Value = '  +2.62%  '

# Row 47: NEARProtocol
$ws.Range('E47').Value = '  +6.22%  '

# Row 48: FraxShare
$ws.Range('D48').Value = '''10.20'
$ws.Range('E48').Value = '  +5.37%  '

# Row 49: MultiversX
$ws.Range('D49').Value = '''57.06'
$ws.Range('E49').Value = '  +6.47%  '

# Row 50: HuobiToken
$ws.Range('E50').Value = '  +2.27%  '

# Row 51: Stacks
$ws.Range('E51').Value = '  +7.81%  '
